# 2024 Maas Katsayilari ile diger maas verileri islendi.
# Applies the "2024 Maas Katsayilari" update to the ozel_hizmet_tazminati sheet:
#   - Updates 77 cells (columns C/D, various rows) with the new 2024 figures.
#   - Moves the active tab from the 1st sheet to the 4th sheet
#     (ozel_hizmet_tazminati), matching the user's last-worked-on view.
#   - Updates the selected cell on ozel_hizmet_tazminati to I119.

$wb = $excel.ActiveWorkbook

$wsOzel = $wb.Worksheets.Item("ozel_hizmet_tazminati")

$updates = @{
    "D13" = 160
    "D14" = 160
    "D15" = 150
    "D16" = 150
    "D17" = 140
    "D18" = 140
    "D19" = 140
    "D20" = 140
    "D21" = 140
    "D22" = 170
    "D23" = 170
    "D24" = 160
    "D25" = 160
    "D26" = 150
    "D27" = 150
    "D28" = 150
    "D29" = 150
    "D97" = 115
    "D98" = 115
    "D99" = 105
    "D100" = 105
    "D101" = 95
    "D102" = 95
    "D103" = 95
    "D104" = 90
    "D105" = 90
    "D106" = 90
    "D107" = 115
    "D108" = 115
    "D109" = 105
    "D110" = 105
    "D111" = 95
    "D112" = 95
    "D113" = 95
    "D114" = 90
    "D115" = 90
    "D116" = 90
    "D117" = 115
    "D118" = 115
    "D119" = 105
    "D120" = 105
    "D121" = 95
    "D122" = 95
    "D123" = 95
    "D124" = 90
    "D125" = 90
    "D126" = 90
    "D127" = 115
    "D128" = 115
    "D129" = 105
    "D130" = 105
    "D131" = 95
    "D132" = 95
    "D133" = 95
    "D134" = 90
    "D135" = 90
    "D136" = 90
    "C137" = 140
    "C138" = 140
    "C139" = 140
    "C140" = 140
    "C141" = 139
    "C142" = 139
    "C143" = 39
    "C144" = 139
    "C145" = 139
    "C146" = 139
    "C147" = 138
    "C148" = 138
    "C149" = 138
    "C150" = 138
    "C151" = 137
    "C152" = 137
    "C153" = 137
    "C154" = 137
    "C155" = 137
    "C156" = 136
}

foreach ($cellRef in $updates.Keys) {
    $wsOzel.Range($cellRef).Value = $updates[$cellRef]
}

# Make ozel_hizmet_tazminati the active sheet/tab (activeTab=3, tabSelected
# moves off aylik_gosterge_puanlari and onto ozel_hizmet_tazminati), and leave
# the selection positioned at I119 as in the source edit.
$wsOzel.Activate() | Out-Null
$wsOzel.Range("I119").Select() | Out-Null

